$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- hotel_info sheet: fill in English_Reviews_num, Local_Rank, Total_Reviews_num ---
Set-TextValue $ws1.Range("G2") "9"
Set-TextValue $ws1.Range("H2") "349"
Set-TextValue $ws1.Range("I2") "10"

# --- review_info sheet: add 5 new review rows ---
# row 2
$ws2.Range("A2").Value = 5218
$ws2.Range("D2").Value = 1
Set-TextValue $ws2.Range("E2") "08/12/2018"
Set-TextValue $ws2.Range("F2") "https://www.tripadvisor.com/ShowUserReviews-g56003-d11904353-r479022334-Red_Roof_Inn_Houston_East-Houston_Texas.html"
Set-TextValue $ws2.Range("G2") "56003"
Set-TextValue $ws2.Range("H2") "11904353"
Set-TextValue $ws2.Range("I2") "479022334"
Set-TextValue $ws2.Range("J2") "04/25/2017"
Set-TextValue $ws2.Range("K2") "business  stay  with this motel"
Set-TextValue $ws2.Range("L2") "i had great experience with staying in this perticular red roof inn. good staff with helpful and kind natured . clean and neat room with good location with restaurant next to it.  breakfast was good . i liked to stay with them. overall very convenient stay with this red roof inn.MoreShow less"
$ws2.Range("M2").Value = 4
Set-TextValue $ws2.Range("N2") "March 2017"
Set-TextValue $ws2.Range("O2") " traveled on business"
$ws2.Range("Q2").Value = 4
$ws2.Range("U2").Value = 4
$ws2.Range("V2").Value = 0
Set-TextValue $ws2.Range("W2") "Texas45Lodging, Guest Relations Manager at Red Roof Inn Houston East, responded to this reviewResponded July 26, 2017"
Set-TextValue $ws2.Range("X2") "Responded July 26, 2017"
Set-TextValue $ws2.Range("Y2") "i had great experience with staying in this perticular red roof inn. good staff with helpful and kind natured . clean and neat room with good location with restaurant next to it.  breakfast was good . i liked to stay with them. overall very convenient stay with this red roof inn.More"

# row 3
$ws2.Range("A3").Value = 5218
$ws2.Range("D3").Value = 2
Set-TextValue $ws2.Range("E3") "08/12/2018"
Set-TextValue $ws2.Range("F3") "https://www.tripadvisor.com/ShowUserReviews-g56003-d11904353-r470591380-Red_Roof_Inn_Houston_East-Houston_Texas.html"
Set-TextValue $ws2.Range("G3") "56003"
Set-TextValue $ws2.Range("H3") "11904353"
Set-TextValue $ws2.Range("I3") "470591380"
Set-TextValue $ws2.Range("J3") "03/27/2017"
Set-TextValue $ws2.Range("K3") "Great Place for the Price!"
Set-TextValue $ws2.Range("L3") "Great place for a great price. Friendly staff and great rooms for the price. Everything's around the corner from food places to movie theater. From checking in to checking out my experience was great. Definitely staying here again.MoreShow less"
$ws2.Range("M3").Value = 5
Set-TextValue $ws2.Range("N3") "March 2017"
Set-TextValue $ws2.Range("O3") " traveled on business"
$ws2.Range("S3").Value = 5
$ws2.Range("U3").Value = 5
$ws2.Range("V3").Value = 0
Set-TextValue $ws2.Range("W3") "Texas45Lodging, Guest Relations Manager at Red Roof Inn Houston East, responded to this reviewResponded July 26, 2017"
Set-TextValue $ws2.Range("X3") "Responded July 26, 2017"
Set-TextValue $ws2.Range("Y3") "Great place for a great price. Friendly staff and great rooms for the price. Everything's around the corner from food places to movie theater. From checking in to checking out my experience was great. Definitely staying here again.More"

# row 4
$ws2.Range("A4").Value = 5218
$ws2.Range("D4").Value = 3
Set-TextValue $ws2.Range("E4") "08/12/2018"
Set-TextValue $ws2.Range("F4") "https://www.tripadvisor.com/ShowUserReviews-g56003-d11904353-r470177325-Red_Roof_Inn_Houston_East-Houston_Texas.html"
Set-TextValue $ws2.Range("G4") "56003"
Set-TextValue $ws2.Range("H4") "11904353"
Set-TextValue $ws2.Range("I4") "470177325"
Set-TextValue $ws2.Range("J4") "03/25/2017"
Set-TextValue $ws2.Range("K4") "If your choice is sleeping in the car or staying here, sleep in the car."
Set-TextValue $ws2.Range("L4") "Probably was a nice place before it filled with construction workers. Seriously run down, stained carpets, curtains falling off, holes in the bedding, and seriously funky odor in the room.  At least the bathroom and towels are clean. Pleasant lobby. And they do indeed allow pets.MoreShow less"
$ws2.Range("M4").Value = 2
Set-TextValue $ws2.Range("N4") "March 2017"
Set-TextValue $ws2.Range("O4") " traveled as a couple"
$ws2.Range("V4").Value = 0
Set-TextValue $ws2.Range("W4") "Texas45Lodging, Guest Relations Manager at Red Roof Inn Houston East, responded to this reviewResponded July 26, 2017"
Set-TextValue $ws2.Range("X4") "Responded July 26, 2017"
Set-TextValue $ws2.Range("Y4") "Probably was a nice place before it filled with construction workers. Seriously run down, stained carpets, curtains falling off, holes in the bedding, and seriously funky odor in the room.  At least the bathroom and towels are clean. Pleasant lobby. And they do indeed allow pets.More"

# row 5
$ws2.Range("A5").Value = 5218
$ws2.Range("D5").Value = 4
Set-TextValue $ws2.Range("E5") "08/12/2018"
Set-TextValue $ws2.Range("F5") "https://www.tripadvisor.com/ShowUserReviews-g56003-d11904353-r461820521-Red_Roof_Inn_Houston_East-Houston_Texas.html"
Set-TextValue $ws2.Range("G5") "56003"
Set-TextValue $ws2.Range("H5") "11904353"
Set-TextValue $ws2.Range("I5") "461820521"
Set-TextValue $ws2.Range("J5") "02/21/2017"
Set-TextValue $ws2.Range("K5") "Disgusting"
Set-TextValue $ws2.Range("L5") "I am a member of the Red Roof Inn rewards program and frequently stay at their hotels. They are never particularly nice but this location has got to be the worst one I've ever stayed at. A majority of the exterior lighting is either missing or doesn't work, the blankets on my beds had burn holes in them, the closet door had a hole in it, there were plumbing leaks in the bathroom, and there was some kind of sewer line leak in the rear parking lot area where there is standing water that had an unpleasant odor.MoreShow less"
$ws2.Range("M5").Value = 1
Set-TextValue $ws2.Range("N5") "February 2017"
Set-TextValue $ws2.Range("O5") " traveled solo"
$ws2.Range("P5").Value = 2
$ws2.Range("Q5").Value = 2
$ws2.Range("U5").Value = 3
$ws2.Range("V5").Value = 0
Set-TextValue $ws2.Range("W5") "Texas45Lodging, Guest Relations Manager at Red Roof Inn Houston East, responded to this reviewResponded July 26, 2017"
Set-TextValue $ws2.Range("X5") "Responded July 26, 2017"
Set-TextValue $ws2.Range("Y5") "I am a member of the Red Roof Inn rewards program and frequently stay at their hotels. They are never particularly nice but this location has got to be the worst one I've ever stayed at. A majority of the exterior lighting is either missing or doesn't work, the blankets on my beds had burn holes in them, the closet door had a hole in it, there were plumbing leaks in the bathroom, and there was some kind of sewer line leak in the rear parking lot area where there is standing water that had an unpleasant odor.More"

# row 6
$ws2.Range("A6").Value = 5218
$ws2.Range("D6").Value = 5
Set-TextValue $ws2.Range("E6") "08/12/2018"
Set-TextValue $ws2.Range("F6") "https://www.tripadvisor.com/ShowUserReviews-g56003-d11904353-r460882525-Red_Roof_Inn_Houston_East-Houston_Texas.html"
Set-TextValue $ws2.Range("G6") "56003"
Set-TextValue $ws2.Range("H6") "11904353"
Set-TextValue $ws2.Range("I6") "460882525"
Set-TextValue $ws2.Range("J6") "02/18/2017"
Set-TextValue $ws2.Range("K6") "Horrible"
Set-TextValue $ws2.Range("L6") "This place is disgusting.  There is mold in the rooms.  The buildings are very old and have obvious poorly done patch jobs.  I have stayed in many Red Roof Inns over the years.  The company should be ashamed of this place.  It should be torn down.MoreShow less"
$ws2.Range("M6").Value = 1
Set-TextValue $ws2.Range("N6") "February 2017"
Set-TextValue $ws2.Range("O6") " traveled solo"
$ws2.Range("R6").Value = 1
$ws2.Range("U6").Value = 1
$ws2.Range("V6").Value = 0
Set-TextValue $ws2.Range("W6") "Texas45Lodging, Guest Relations Manager at Red Roof Inn Houston East, responded to this reviewResponded July 26, 2017"
Set-TextValue $ws2.Range("X6") "Responded July 26, 2017"
Set-TextValue $ws2.Range("Y6") "This place is disgusting.  There is mold in the rooms.  The buildings are very old and have obvious poorly done patch jobs.  I have stayed in many Red Roof Inns over the years.  The company should be ashamed of this place.  It should be torn down.More"
